# Actualizacion automatica de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversion del dia" note with new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldText = $ws1.Range("A1").Value()
$newText = $oldText.Replace("✅ 1000 Bs = 7.84 = 31857.69 pesos", "✅ 1000 Bs = 7.98 = 32268.13 pesos")
$newText = $newText.Replace("✅ 31857.69 pesos = 7.83 = 951.58 Bs", "✅ 32268.13 pesos = 7.94 = 944.99 Bs")
$ws1.Range("A1").Value = $newText

# --- tasas: update the N10/O10/N12/O12 rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 125.359
$ws2.Range("O10").Value = 4045.1
$ws2.Range("N12").Value = 4063.44
$ws2.Range("O12").Value = 119
